$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns F, G, H
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Apply the same style as the existing header row (A1:E1) to the new headers
$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats

# Fill boolean FALSE values for rows 2 through 12 in columns F, G, H
for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 6).Value = $false
    $ws.Cells.Item($row, 7).Value = $false
    $ws.Cells.Item($row, 8).Value = $false
}
